$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.319.13"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.12%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.671.52"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.46%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'685.28"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.21%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'159.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.20%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -1.44%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.25%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.77%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.434"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -3.46%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.0000232"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.84%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.291.90"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.40%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'32.17"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.94%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = "'WrappedBTC"
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = "'69.327.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.02%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = "'WrappedEther"
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = "'3.661.97"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.76%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.04%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'15.78"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -3.20%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.37"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.81%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'468.41"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.35%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.647"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.63%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'79.66"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.35%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'3.818.65"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.37%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.04%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.0000122"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -4.41%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'10.89"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -5.39%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'9.18"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -4.01%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.77%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.73"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -6.00%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'6.56"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -3.11%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.08%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -6.05%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'26.78"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.67%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'3.645.28"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.25%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.159"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.20%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'8.14"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -4.73%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'6.12"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.10%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D40').Value = "'2.21"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.06%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.0895"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.05%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.01%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'Monero"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'166.51"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +5.86%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Mantle"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.939"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -2.10%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.45%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.000281"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.16%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'SuiNetwork"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'1.11"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +4.83%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'dogwifhat"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'2.70"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.76%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.44%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'27.32"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.84%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -3.76%  "
$ws.Range('E51').Style = 'Normal'
